# Aspekt 3 im Projekthandbuch
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 52 originally held 13.07.2013 (serial 41468) - it actually belongs to
# the 14.07.2013 entry (serial 41469).
$ws.Range("A52").Value = 41469

# New row 53: 15.07.2013, Handbuch Aspekt 3 work, 2.5 hours.
$ws.Range("A53").Value = 41470
$ws.Range("B53").Value = "Handbuch, Aspekt 3 Commands / Handler, Menus, Key Bindings"
$ws.Range("C53").Value = 2.5

# TOTAL row: extend the sum to include the newly added row 53.
$ws.Range("C57").Formula = "=SUM(C2:C53)"

# Reflect the scrolled/selected view from the edit session.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("C53").Select()
